# Append a new paragraph after the "Resources" heading paragraph containing:
#   "Design idea for pet card: https://www.dspca.ie/adoption-cats-landing-page/"
# as two separate runs, matching the target OOXML exactly (no inherited
# Heading1 style, no leftover empty paragraph, original paragraph untouched).

$d = $word.ActiveDocument

# Locate the "Resources" paragraph via Find so the script is resilient even
# if the paragraph's position in the document changes.
$searchRange = $d.Content
$found = $searchRange.Find.Execute("Resources", $false, $false, $false, $false,
                                    $false, $true, 1, $false, "", 0)

if ($found -and $searchRange.Text -eq "Resources") {
    $targetParaRange = $searchRange.Paragraphs(1).Range
} else {
    # Fallback: use the first paragraph of the document.
    $targetParaRange = $d.Paragraphs(1).Range
}

# Position right before that paragraph's trailing paragraph mark, so the
# inserted paragraph is spliced in immediately after it without disturbing
# the existing paragraph or leaving a stray empty paragraph behind.
$insertPos = $targetParaRange.End - 1
$insertionPoint = $d.Range($insertPos, $insertPos)

$newParagraphXml = '<?xml version="1.0" standalone="yes"?>' +
'<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
'<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
'<pkg:xmlData>' +
'<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
'<w:body>' +
'<w:p>' +
'<w:r><w:t xml:space="preserve">Design idea for pet card: </w:t></w:r>' +
'<w:r><w:t>https://www.dspca.ie/adoption-cats-landing-page/</w:t></w:r>' +
'</w:p>' +
'</w:body>' +
'</w:document>' +
'</pkg:xmlData>' +
'</pkg:part>' +
'</pkg:package>'

[void]$insertionPoint.InsertXML($newParagraphXml)
